$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 94) that continues the existing monthly data series.
# Copy the format from the row above (A93:C93) down to row 94 first so that
# the date cell (column A) keeps the same date number format/style as the
# rest of the column, then overwrite the values.
$ws.Range("A93:C93").Copy($ws.Range("A94:C94"))

$ws.Range("A94").Value = 45536
$ws.Range("B94").Value = 0.135766982849478
$ws.Range("C94").Value = 0.0528682496322757
